# Daily update at 8 AM UTC
# Appends the next day's row to the "Wins Over Time" log and moves the
# "latest row" date formatting down to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 77
$newRow = $lastRow + 1

# The previous "last row" (A77) used a distinct date-only number format to
# highlight the most recent entry. Restore the normal timestamp format
# (matching the rows above it) before adding the new last row.
$ws.Range("A$lastRow").NumberFormat = $ws.Range("A$($lastRow - 1)").NumberFormat

# Append the new day's data.
$ws.Range("A$newRow").Value = 45818
$ws.Range("B$newRow").Value = 334
$ws.Range("C$newRow").Value = 333
$ws.Range("D$newRow").Value = 335

# Give the new last row the distinguishing date-only format.
$ws.Range("A$newRow").NumberFormat = "YYYY-MM-DD"
